# Apply updated cryptocurrency price/volume data to columns D (Price) and E (Volume(1h)).
# Column D values are forced to Text format before assignment (then the style is reset
# back to Normal/default) so that numeric-looking strings such as "303.50" or "1.001"
# are preserved verbatim as text instead of being auto-converted to numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.213.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.73%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.601.58'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3763'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.02'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3632'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.270'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08139'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.73'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.47%  '
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("E15").Value = '  +0.88%  '
$ws.Range("E16").Value = '  +0.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.599.50'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.13'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06920'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.09'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.532'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("E23").Value = '  -1.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.207.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.454'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.050'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '149.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.270'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.89'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.377'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.712'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.777.30'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9643'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07473'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.32'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02737'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.25%  '
$ws.Range("E38").Value = '  -0.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.113'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08773'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.387'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7072'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.38'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.47'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6524'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("E46").Value = '  +0.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.312'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.15%  '
$ws.Range("E48").Value = '  +0.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '132.11'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07926'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.203'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.91%  '
